$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$V$13') {
        $h.Address = "mailto:pqrs321@abc.com"
        $h.TextToDisplay = "pqrs321@abc.com"
    }
}
foreach ($h in $ws.Hyperlinks) {
    Write-Host $h.Address
    Write-Host $h.TextToDisplay
}
